$d = $word.ActiveDocument

function Find-ParaIndex($pattern) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# Locate the Q.13 question paragraph; the blank paragraph right before it is
# the one to remove, and two paragraphs before it is the "(Answer found
# here)" paragraph whose explicit paragraph formatting needs to go away.
$q13Idx = Find-ParaIndex("Q\.13 \[line 160\]")
$blankIdx = $q13Idx - 1
$answerIdx = $q13Idx - 2

# 1) Delete the blank paragraph sitting between the "(Answer found here)"
#    paragraph and the Q.13 paragraph.
$d.Paragraphs.Item($blankIdx).Range.Delete()

# 2) Strip the paragraph formatting (the explicit "spacing after 0") from
#    the "(Answer found here)" paragraph, leaving it with no <w:pPr/>, while
#    keeping its text/hyperlink content intact.
$answerPara = $d.Paragraphs.Item($answerIdx)
$answerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:r><w:t xml:space="preserve">(Answer found here) </w:t></w:r><w:hyperlink r:id="rId13" w:history="1"><w:r><w:t>http://www.cplusplus.com/forum/beginner/138279/</w:t></w:r></w:hyperlink></w:p>'
$answerPara.Range.InsertXML($answerXml) | Out-Null
# Re-apply the Hyperlink character style (InsertXML drops rStyle refs), so
# the link keeps its original look.
$answerPara2 = $d.Paragraphs.Item($answerIdx)
$answerPara2.Range.Hyperlinks.Item(1).Range.Style = "Hyperlink"

# 3) Add a new answer paragraph "p1" right after the Q.13 question.
$q13Idx2 = Find-ParaIndex("Q\.13 \[line 160\]")
$q13Para = $d.Paragraphs.Item($q13Idx2)
$q13Para.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($q13Idx2 + 1)
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>p1</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newXml) | Out-Null
